$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 24.999262
$ws.Range("H2").Value = 74.997786
$ws.Range("I2").Value = 0.2094245171924971
$ws.Range("J2").Value = 0.209424517192497
$ws.Range("M2").Value = 8.533046666666666
$ws.Range("N2").Value = 25.59914
$ws.Range("O2").Value = 0.2932132236642383
$ws.Range("P2").Value = 0.2932132236642383
$ws.Range("Q2").Value = 213.3198692782267
$ws.Range("R2").Value = 1919.87882350404
$ws.Range("S2").Value = 0.06140603780033876
$ws.Range("T2").Value = 0.06140603780033875

# Row 3
$ws.Range("G3").Value = 24.999262
$ws.Range("H3").Value = 74.997786
$ws.Range("I3").Value = 0.2094245171924971
$ws.Range("J3").Value = 0.209424517192497
$ws.Range("O3").Value = 0.3119288965200195
$ws.Range("P3").Value = 0.3119288965200194
$ws.Range("Q3").Value = 226.9359839853214
$ws.Range("R3").Value = 2042.423855867892
$ws.Range("S3").Value = 0.06532555855209346
$ws.Range("T3").Value = 0.06532555855209343

# Row 4
$ws.Range("G4").Value = 24.999262
$ws.Range("H4").Value = 74.997786
$ws.Range("I4").Value = 0.2094245171924971
$ws.Range("J4").Value = 0.209424517192497
$ws.Range("O4").Value = 0.3948578798157423
$ws.Range("P4").Value = 0.3948578798157423
$ws.Range("Q4").Value = 287.268869572628
$ws.Range("R4").Value = 2585.419826153652
$ws.Range("S4").Value = 0.08269292084006485
$ws.Range("T4").Value = 0.08269292084006484

# Row 5
$ws.Range("I5").Value = 0.522807373179233
$ws.Range("J5").Value = 0.5228073731792329
$ws.Range("M5").Value = 8.533046666666666
$ws.Range("N5").Value = 25.59914
$ws.Range("O5").Value = 0.2932132236642383
$ws.Range("P5").Value = 0.2932132236642383
$ws.Range("Q5").Value = 532.5317302834044
$ws.Range("R5").Value = 4792.78557255064
$ws.Range("S5").Value = 0.1532940352453153
$ws.Range("T5").Value = 0.1532940352453153

# Row 6
$ws.Range("I6").Value = 0.522807373179233
$ws.Range("J6").Value = 0.5228073731792329
$ws.Range("O6").Value = 0.3119288965200195
$ws.Range("P6").Value = 0.3119288965200194
$ws.Range("S6").Value = 0.1630787270083282
$ws.Range("T6").Value = 0.1630787270083281

# Row 7
$ws.Range("I7").Value = 0.522807373179233
$ws.Range("J7").Value = 0.5228073731792329
$ws.Range("O7").Value = 0.3948578798157423
$ws.Range("P7").Value = 0.3948578798157423
$ws.Range("S7").Value = 0.2064346109255895
$ws.Range("T7").Value = 0.2064346109255895

# Row 8
$ws.Range("I8").Value = 0.2677681096282701
$ws.Range("J8").Value = 0.2677681096282701
$ws.Range("M8").Value = 8.533046666666666
$ws.Range("N8").Value = 25.59914
$ws.Range("O8").Value = 0.2932132236642383
$ws.Range("P8").Value = 0.2932132236642383
$ws.Range("Q8").Value = 272.7486681527222
$ws.Range("R8").Value = 2454.7380133745
$ws.Range("S8").Value = 0.07851315061858423
$ws.Range("T8").Value = 0.07851315061858423

# Row 9
$ws.Range("I9").Value = 0.2677681096282701
$ws.Range("J9").Value = 0.2677681096282701
$ws.Range("O9").Value = 0.3119288965200195
$ws.Range("P9").Value = 0.3119288965200194
$ws.Range("S9").Value = 0.08352461095959789
$ws.Range("T9").Value = 0.08352461095959787

# Row 10
$ws.Range("I10").Value = 0.2677681096282701
$ws.Range("J10").Value = 0.2677681096282701
$ws.Range("O10").Value = 0.3948578798157423
$ws.Range("P10").Value = 0.3948578798157423
$ws.Range("S10").Value = 0.105730348050088
$ws.Range("T10").Value = 0.105730348050088
